$wb = $excel.ActiveWorkbook

# Deselect the current range on EmpList and select A1:B1 instead,
# matching the state captured once the new sheet becomes active.
$empList = $wb.Worksheets.Item("EmpList")
$empList.Range("A1:B1").Select()

# Add the new "ApplyLeave" sheet right after the last existing sheet
# (EmpList), which also makes it the new active sheet/tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$applyLeave = $wb.Worksheets.Add($null, $lastSheet)
$applyLeave.Name = "ApplyLeave"

# Populate the header row for the new test-case data sheet.
$applyLeave.Range("A1").Value = "Admin"
$applyLeave.Range("B1").Value = "admin123"
$applyLeave.Range("C1").Value = "By.xpath(`"//*[text()='Leave Type ']//following::select[1]`")"
$applyLeave.Range("D1").Value = "By.xpath(`"//table[contains(@class,'calendar')]//a[text()=25]`")"
$applyLeave.Range("E1").Value = "I am out sick"

# Leave the selection positioned at D16 on the new sheet.
$applyLeave.Range("D16").Select()
